# Sprint and Product Backlog 1
# Fill in the Sprint-1 "Logs" burndown/story-point tracker with the actual
# per-day story point values for each team member, and fix the duplicated
# date in the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$ws.Activate()

# Row 4 had a duplicated date (45014 appeared twice in F4/G4) - the F4 date
# should actually be 45013.
$ws.Range("F4").Value = 45013

# Story point entries for Thu/Fri/Mon (columns C/D/E) for each team member.
$values = @(
    @(5.5, 6,   5),
    @(6,   6,   5.5),
    @(6,   6.5, 5),
    @(5,   5.5, 6),
    @(5,   6,   5.5),
    @(5,   5,   6),
    @(5.5, 6,   5.5),
    @(6,   7,   8),
    @(5.5, 7,   8),
    @(5,   5,   6)
)

$row = 6
foreach ($entry in $values) {
    $ws.Range("C$row").Value = $entry[0]
    $ws.Range("D$row").Value = $entry[1]
    $ws.Range("E$row").Value = $entry[2]
    $row++
}

# Center-align the whole tracker grid (values + still-blank cells).
$ws.Range("C6:H15").HorizontalAlignment = -4108

$ws.Range("E15").Select()
